# Applies scheduled market-price/profit data refresh to the Yojimbo Profits workbook.
# For each affected Leve row, updates price/profit columns (H-N) with freshly scraped values.
$wb = $excel.ActiveWorkbook

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1687.3636
$ws.Range("I51").Value = 1450.25
$ws.Range("J51").Value = 1822.8572
$ws.Range("K51").Value = 1450.25
$ws.Range("L51").Value = 1822.8572
$ws.Range("M51").Value = -966.25
$ws.Range("N51").Value = -2790.8572

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2244.1875
$ws.Range("J113").Value = 2715.4546
$ws.Range("L113").Value = 2715.4546
$ws.Range("N113").Value = -9223.454600000001

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 934.82355
$ws.Range("I2").Value = 745.53845
$ws.Range("J2").Value = 1550
$ws.Range("K2").Value = 745.53845
$ws.Range("L2").Value = 1550
$ws.Range("M2").Value = -632.53845
$ws.Range("N2").Value = -1776

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1313
$ws.Range("I45").Value = 1313
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1313
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -936
$ws.Range("N45").ClearContents()

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 934.82355
$ws.Range("I116").Value = 745.53845
$ws.Range("J116").Value = 1550
$ws.Range("K116").Value = 745.53845
$ws.Range("L116").Value = 1550
$ws.Range("M116").Value = 1548.46155
$ws.Range("N116").Value = -6138

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 43653.668
$ws.Range("J133").Value = 43653.668
$ws.Range("L133").Value = 43653.668
$ws.Range("N133").Value = -48713.668

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 934.82355
$ws.Range("I3").Value = 745.53845
$ws.Range("J3").Value = 1550
$ws.Range("K3").Value = 745.53845
$ws.Range("L3").Value = 1550
$ws.Range("M3").Value = -631.53845
$ws.Range("N3").Value = -1778

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1494.091
$ws.Range("I99").Value = 1433.5714
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 1433.5714
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = 64.42859999999996
$ws.Range("N99").Value = -4596

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 939
$ws.Range("I16").Value = 925
$ws.Range("J16").Value = 995
$ws.Range("K16").Value = 925
$ws.Range("L16").Value = 995
$ws.Range("M16").Value = -638
$ws.Range("N16").Value = -1569

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 797.2273
$ws.Range("I22").Value = 920.73334
$ws.Range("J22").Value = 532.5714
$ws.Range("K22").Value = 920.73334
$ws.Range("L22").Value = 532.5714
$ws.Range("M22").Value = -570.73334
$ws.Range("N22").Value = -1232.5714

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2090.5557
$ws.Range("I62").Value = 2069.1667
$ws.Range("J62").Value = 2133.3333
$ws.Range("K62").Value = 2069.1667
$ws.Range("L62").Value = 2133.3333
$ws.Range("M62").Value = -1445.1667
$ws.Range("N62").Value = -3381.3333

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2090.5557
$ws.Range("I65").Value = 2069.1667
$ws.Range("J65").Value = 2133.3333
$ws.Range("K65").Value = 10345.8335
$ws.Range("L65").Value = 10666.6665
$ws.Range("M65").Value = -7225.833500000001
$ws.Range("N65").Value = -16906.6665

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 810
$ws.Range("I105").Value = 650
$ws.Range("J105").Value = 1050
$ws.Range("K105").Value = 650
$ws.Range("L105").Value = 1050
$ws.Range("M105").Value = 1097
$ws.Range("N105").Value = -4544

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 939
$ws.Range("I113").Value = 925
$ws.Range("J113").Value = 995
$ws.Range("K113").Value = 925
$ws.Range("L113").Value = 995
$ws.Range("M113").Value = 1245
$ws.Range("N113").Value = -5335

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6960.2666
$ws.Range("I68").Value = 316.83334
$ws.Range("J68").Value = 33534
$ws.Range("K68").Value = 950.5000200000001
$ws.Range("L68").Value = 100602
$ws.Range("M68").Value = -139.5000200000001
$ws.Range("N68").Value = -102224

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 6960.2666
$ws.Range("I71").Value = 316.83334
$ws.Range("J71").Value = 33534
$ws.Range("K71").Value = 2851.50006
$ws.Range("L71").Value = 301806
$ws.Range("M71").Value = 1204.49994
$ws.Range("N71").Value = -309918

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1365.1333
$ws.Range("I121").Value = 479.7
$ws.Range("J121").Value = 3136
$ws.Range("K121").Value = 1439.1
$ws.Range("L121").Value = 9408
$ws.Range("M121").Value = -129.0999999999999
$ws.Range("N121").Value = -12028

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2488.111
$ws.Range("I140").Value = 1940.8182
$ws.Range("J140").Value = 2864.375
$ws.Range("K140").Value = 5822.4546
$ws.Range("L140").Value = 8593.125
$ws.Range("M140").Value = -642.4546
$ws.Range("N140").Value = -18953.125

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3938.0715
$ws.Range("I70").Value = 3921.138
$ws.Range("J70").Value = 3975.8462
$ws.Range("K70").Value = 3921.138
$ws.Range("L70").Value = 3975.8462
$ws.Range("M70").Value = -3651.138
$ws.Range("N70").Value = -4515.8462

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 3938.0715
$ws.Range("I73").Value = 3921.138
$ws.Range("J73").Value = 3975.8462
$ws.Range("K73").Value = 3921.138
$ws.Range("L73").Value = 3975.8462
$ws.Range("M73").Value = -2985.138
$ws.Range("N73").Value = -5847.8462

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1301.5
$ws.Range("I126").Value = 1167.174
$ws.Range("J126").Value = 1742.8572
$ws.Range("K126").Value = 3501.522
$ws.Range("L126").Value = 5228.571599999999
$ws.Range("M126").Value = -1031.522
$ws.Range("N126").Value = -10168.5716

# GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 13850.667
$ws.Range("J134").Value = 13850.667
$ws.Range("L134").Value = 41552.001
$ws.Range("N134").Value = -46622.001

# GSM row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 57757.145
$ws.Range("J135").Value = 57757.145
$ws.Range("L135").Value = 57757.145
$ws.Range("N135").Value = -67897.14499999999

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2038.6875
$ws.Range("I61").Value = 1601.9
$ws.Range("K61").Value = 1601.9
$ws.Range("M61").Value = -1399.9

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 879.63635
$ws.Range("I93").Value = 923.25
$ws.Range("J93").Value = 763.3333
$ws.Range("K93").Value = 923.25
$ws.Range("L93").Value = 763.3333
$ws.Range("M93").Value = 324.75
$ws.Range("N93").Value = -3259.3333

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2038.6875
$ws.Range("I113").Value = 1601.9
$ws.Range("K113").Value = 1601.9
$ws.Range("M113").Value = 568.0999999999999
